# ASNN-335 - Design for instructor preview of assignment
# Insert a new TOC row ("PREVIEW AS STUDENT" -> page 25) right before the
# "EMAIL NOTIFICATION FOR STUDENT SUBMISSION" row, and bump every
# subsequent page number in the Table of Contents by 4.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Helper: replace a cell's paragraph contents with one or more runs that
# all share the page-number run formatting used throughout this table
# (Trebuchet MS, bold, sz 24 / szCs 24). Passing more than one string in
# $parts produces that many separate <w:r> runs (used where the source
# diff shows the text typed across two runs).
# ---------------------------------------------------------------------
function Set-PageNumberCell($cell, [string[]]$parts) {
    $runsXml = ""
    foreach ($part in $parts) {
        $runsXml += '<w:r ' + $wNs + '><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + $part + '</w:t></w:r>'
    }
    $paraXml = '<w:p ' + $wNs + '><w:pPr><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + $runsXml + '</w:p>'
    $para = $cell.Range.Paragraphs.Item(1)
    $para.Range.InsertXML($paraXml)
}

# ---------------------------------------------------------------------
# Step 1: locate the "EMAIL NOTIFICATION FOR STUDENT SUBMISSION" row and
# insert a brand-new row above it, cloning its formatting.
# ---------------------------------------------------------------------
$emailRowIndex = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $headText = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($headText -like "*EMAIL NOTIFICATION FOR STUDENT SUBMISSION*") {
        $emailRowIndex = $i
        break
    }
}

$newRow = $t.Rows.Add($t.Rows.Item($emailRowIndex))
$newRow.Cells.Item(1).Range.Text = "PREVIEW AS STUDENT"
Set-PageNumberCell $newRow.Cells.Item(2) @("25")

# ---------------------------------------------------------------------
# Step 2: walk the remaining rows (the inserted row pushed everything
# down by one) and rewrite each page-number cell per the diff.
# ---------------------------------------------------------------------
$updates = @(
    @{ Title = "EMAIL NOTIFICATION FOR STUDENT SUBMISSION";        Parts = @("29") },
    @{ Title = "ADD GRADEBOOK ITEM HELPER";                        Parts = @("31") },
    @{ Title = "ADD ATTACHMENT HELPER";                            Parts = @("32") },
    @{ Title = "EDIT ASSIGNMENT";                                  Parts = @("33") },
    @{ Title = "EDIT ASSIGNMENT: SUBMISSIONS NOT REQUIRED";        Parts = @("37") },
    @{ Title = "VIEW GRADED SUBMISSIONS";                          Parts = @("38") },
    @{ Title = "RELEASE GRADES HELPER";                            Parts = @("45") },
    @{ Title = "RELEASE FEEDBACK HELPER";                          Parts = @("46") },
    @{ Title = "VIEW UNGRADED SUBMISSIONS";                        Parts = @("4", "8") },
    @{ Title = [string]::Concat("VIEW STUDENT", [char]0x2019, "S SUBMISSION"); Parts = @("50") },
    @{ Title = "PROVIDE FEEDBACK FOR NE";                          Parts = @("5", "8") },
    @{ Title = "DOWNLOAD ALL";                                     Parts = @("60") },
    @{ Title = "UPLOAD";                                           Parts = @("63") },
    @{ Title = "SV: ASSIGNMENT LIST WITH NO ASSIGNMENTS";          Parts = @("6", "5") },
    @{ Title = "SV: ASSIGNMENT LIST";                              Parts = @("6", "6"); Exact = $true },
    @{ Title = "SETTINGS FOR STUDENT";                             Parts = @("70") },
    @{ Title = "EMAIL CONFIRMATION FOR SUCCESSFUL SUBMISSION";     Parts = @("72") },
    @{ Title = "SV: VIEW DETAILS AND SUBMIT";                      Parts = @("73") },
    @{ Title = "SV: PREVIEW DETAILS AND SUBMIT";                   Parts = @("81") },
    @{ Title = "SV: SUBMISSION CONFIRMATION";                      Parts = @("83") },
    @{ Title = "SV: ASSIGNMENT LIST CONFIRMATION";                 Parts = @("8", "6") },
    @{ Title = "SV: VIEW SUBMISSION";                              Parts = @("8", "8") },
    @{ Title = "SV: MULTIPLE SUBMISSIONS";                         Parts = @("91") },
    @{ Title = "SV: RESUBMIT";                                     Parts = @("9", "5") },
    @{ Title = "SV: VIEW DETAILS";                                 Parts = @("10", "4"); Exact = $true },
    @{ Title = "AL: DELETED ASSIGNMENT";                           Parts = @("10", "7") },
    @{ Title = "VDS: DELETED ASSIGNMENT WITH SAVED DRAFT";         Parts = @("10", "9") },
    @{ Title = "SVVS: DELETED ASSIGNMENT WITH SUBMISSION";         Parts = @("1", "13") },
    @{ Title = "LINK FROM SCHEDULE PERMISSIONS ERROR";             Parts = @("11", "7") }
)

$searchStart = $emailRowIndex
for ($u = 0; $u -lt $updates.Count; $u++) {
    $entry = $updates[$u]
    for ($i = $searchStart; $i -le $t.Rows.Count; $i++) {
        $headText = $t.Rows.Item($i).Cells.Item(1).Range.Text
        $isMatch = $false
        if ($entry.Exact) {
            $isMatch = ($headText -eq ($entry.Title + "`r`a")) -or ($headText -eq $entry.Title)
        } else {
            $isMatch = $headText -like ("*" + $entry.Title + "*")
        }
        if ($isMatch) {
            Set-PageNumberCell $t.Rows.Item($i).Cells.Item(2) $entry.Parts
            $searchStart = $i + 1
            break
        }
    }
}
